$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new I and J columns - copy formatting from H1 (bold, border, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF values per row (row number => [I, J])
$values = @{
    2  = @(1, 6)
    3  = @(1, 4)
    4  = @(1, 6)
    5  = @(1, 5)
    6  = @(1, 3)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 6)
    10 = @(3, 8)
    11 = @(1, 4)
    12 = @(1, 4)
    13 = @(1, 6)
    14 = @(1, 7)
    15 = @(1, 6)
    16 = @(1, 5)
    17 = @(1, 7)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 3)
    21 = @(1, 4)
    22 = @(1, 6)
    23 = @(1, 5)
    24 = @(1, 5)
    25 = @(1, 5)
    26 = @(6, 9)
    27 = @(1, 4)
    28 = @(1, 3)
    29 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
